$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 39.6
$ws.Range("I6").Value = 24.5
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 73.5
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = 38.5
$ws.Range("N6").Value = -524

$ws.Range("H8").Value = 2721.875
$ws.Range("I8").Value = 2721.875
$ws.Range("K8").Value = 8165.625
$ws.Range("M8").Value = -8026.625

$ws.Range("H21").Value = 3200
$ws.Range("J21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("N21").Value = -2436

$ws.Range("H23").Value = 3200
$ws.Range("J23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1968

$ws.Range("H31").Value = 628.8889
$ws.Range("I31").Value = 276.66666
$ws.Range("J31").Value = 1333.3334
$ws.Range("K31").Value = 829.9999799999999
$ws.Range("L31").Value = 4000.0002
$ws.Range("M31").Value = -599.9999799999999
$ws.Range("N31").Value = -4460.0002

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H127").Value = 903.125
$ws.Range("I127").Value = 889.2857
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 2667.8571
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 2292.1429
$ws.Range("N127").Value = -12920

$ws.Range("H137").Value = 5141.147
$ws.Range("I137").Value = 1990.6818
$ws.Range("J137").Value = 10917
$ws.Range("K137").Value = 5972.0454
$ws.Range("L137").Value = 32751
$ws.Range("M137").Value = -3422.0454
$ws.Range("N137").Value = -37851

$ws.Range("H138").Value = 3339.3845
$ws.Range("J138").Value = 3834.64
$ws.Range("L138").Value = 11503.92
$ws.Range("N138").Value = -21783.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9606.272000000001
$ws.Range("I45").Value = 11297.444
$ws.Range("K45").Value = 11297.444
$ws.Range("M45").Value = -10920.444

$ws.Range("H74").Value = 205686.06
$ws.Range("I74").Value = 326777.53
$ws.Range("K74").Value = 326777.53
$ws.Range("M74").Value = -325903.53

$ws.Range("H77").Value = 205686.06
$ws.Range("I77").Value = 326777.53
$ws.Range("K77").Value = 1633887.65
$ws.Range("M77").Value = -1629519.65

$ws.Range("H102").Value = 1496.1
$ws.Range("I102").Value = 1496.1
$ws.Range("K102").Value = 1496.1
$ws.Range("M102").Value = 125.9000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2805.325
$ws.Range("I31").Value = 1168.75
$ws.Range("J31").Value = 5260.1875
$ws.Range("K31").Value = 1168.75
$ws.Range("L31").Value = 5260.1875
$ws.Range("M31").Value = -873.75
$ws.Range("N31").Value = -5850.1875

$ws.Range("H34").Value = 2805.325
$ws.Range("I34").Value = 1168.75
$ws.Range("J34").Value = 5260.1875
$ws.Range("K34").Value = 1168.75
$ws.Range("L34").Value = 5260.1875
$ws.Range("M34").Value = -966.75
$ws.Range("N34").Value = -5664.1875

$ws.Range("H58").Value = 3755.7188
$ws.Range("I58").Value = 1531.7894
$ws.Range("J58").Value = 7006.077
$ws.Range("K58").Value = 1531.7894
$ws.Range("L58").Value = 7006.077
$ws.Range("M58").Value = -1328.7894
$ws.Range("N58").Value = -7412.077

$ws.Range("H99").Value = 3798.4443
$ws.Range("I99").Value = 3731.8333
$ws.Range("K99").Value = 3731.8333
$ws.Range("M99").Value = -2233.8333

$ws.Range("H126").Value = 3798.4443
$ws.Range("I126").Value = 3731.8333
$ws.Range("K126").Value = 11195.4999
$ws.Range("M126").Value = -8725.499899999999

$ws.Range("H136").Value = 3755.7188
$ws.Range("I136").Value = 1531.7894
$ws.Range("J136").Value = 7006.077
$ws.Range("K136").Value = 4595.3682
$ws.Range("L136").Value = 21018.231
$ws.Range("M136").Value = -2045.3682
$ws.Range("N136").Value = -26118.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5649
$ws.Range("J9").Value = 6623.75
$ws.Range("L9").Value = 19871.25
$ws.Range("N9").Value = -20319.25

$ws.Range("H60").Value = 415.15384
$ws.Range("I60").Value = 449.7
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 1349.1
$ws.Range("L60").Value = 900
$ws.Range("M60").Value = -1098.1
$ws.Range("N60").Value = -1402

$ws.Range("H62").Value = 1914.9608
$ws.Range("I62").Value = 1270.7188
$ws.Range("K62").Value = 3812.1564
$ws.Range("M62").Value = -3126.1564

$ws.Range("H65").Value = 1914.9608
$ws.Range("I65").Value = 1270.7188
$ws.Range("K65").Value = 11436.4692
$ws.Range("M65").Value = -8004.469200000001

$ws.Range("H104").Value = 477.8
$ws.Range("I104").Value = 477.8
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 1433.4
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 1187.6
$ws.Range("N104").ClearContents()

$ws.Range("H121").Value = 37037920
$ws.Range("J121").Value = 55556492
$ws.Range("L121").Value = 166669476
$ws.Range("N121").Value = -166672096

$ws.Range("H132").Value = 3579.0952
$ws.Range("I132").Value = 2833.7273
$ws.Range("J132").Value = 4399
$ws.Range("K132").Value = 25503.5457
$ws.Range("L132").Value = 39591
$ws.Range("M132").Value = -22973.5457
$ws.Range("N132").Value = -44651

$ws.Range("H139").Value = 45457580
$ws.Range("I139").Value = 83336060
$ws.Range("J139").Value = 3399.8
$ws.Range("K139").Value = 250008180
$ws.Range("L139").Value = 10199.4
$ws.Range("M139").Value = -250003040
$ws.Range("N139").Value = -20479.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 743851.3
$ws.Range("J26").Value = 34493.168
$ws.Range("L26").Value = 34493.168
$ws.Range("N26").Value = -35053.168

$ws.Range("H50").Value = 743851.3
$ws.Range("J50").Value = 34493.168
$ws.Range("L50").Value = 34493.168
$ws.Range("N50").Value = -35489.168

$ws.Range("H70").Value = 8578.6
$ws.Range("I70").Value = 4079.2
$ws.Range("J70").Value = 13078
$ws.Range("K70").Value = 4079.2
$ws.Range("L70").Value = 13078
$ws.Range("M70").Value = -3809.2
$ws.Range("N70").Value = -13618

$ws.Range("H73").Value = 8578.6
$ws.Range("I73").Value = 4079.2
$ws.Range("J73").Value = 13078
$ws.Range("K73").Value = 4079.2
$ws.Range("L73").Value = 13078
$ws.Range("M73").Value = -3143.2
$ws.Range("N73").Value = -14950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4567.357
$ws.Range("I7").Value = 4575.727
$ws.Range("J7").Value = 4536.6665
$ws.Range("K7").Value = 4575.727
$ws.Range("L7").Value = 4536.6665
$ws.Range("M7").Value = -4463.727
$ws.Range("N7").Value = -4760.6665

$ws.Range("H94").Value = 68250
$ws.Range("J94").Value = 68250
$ws.Range("L94").Value = 68250
$ws.Range("N94").Value = -69602

$ws.Range("H126").Value = 4567.357
$ws.Range("I126").Value = 4575.727
$ws.Range("J126").Value = 4536.6665
$ws.Range("K126").Value = 13727.181
$ws.Range("L126").Value = 13609.9995
$ws.Range("M126").Value = -11257.181
$ws.Range("N126").Value = -18549.9995

$ws.Range("H136").Value = 4529.091
$ws.Range("I136").Value = 2804.5833
$ws.Range("K136").Value = 8413.749899999999
$ws.Range("M136").Value = -5863.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2871.7446
$ws.Range("I132").Value = 1846.8
$ws.Range("J132").Value = 8728.571
$ws.Range("K132").Value = 5540.4
$ws.Range("L132").Value = 26185.713
$ws.Range("M132").Value = -3010.4
$ws.Range("N132").Value = -31245.713

$ws.Range("H136").Value = 1551.9688
$ws.Range("I136").Value = 855.43335
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 2566.30005
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -16.30004999999983
$ws.Range("N136").Value = -41100
